$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.286.70'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.50%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.88'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.37%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.45'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.67%  '

# Row 6
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4671'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.49%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2835'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06521'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.29%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.77'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +8.55%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07933'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.86%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.37'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.55%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.867.21'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.55%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.157'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.40%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6777'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.99'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.64%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.279.44'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.39%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.45'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +7.12%  '

# Row 19
$ws.Range("E19").Value = '  +0.04%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.382'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.44%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.112.24'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007308'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.35%  '

# Row 23
$ws.Range("E23").Value = '  +0.07%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.154'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.25%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.23'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.159'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.46%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.09'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.68%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.930'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.62%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.386'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09720'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.99%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.394'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.476'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.068'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.36%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04738'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.59%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.128'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.59%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7048'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.40%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.33%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01860'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.43%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.574'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.41%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.304'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.50%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.64'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.78%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.958'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8500'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.10%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4180'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.70%  '

# Row 45
$ws.Range("E45").Value = '  +0.03%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.33'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.62%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '965.13'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.36%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.175'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.42%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.295'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.76%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.07'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.03%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1133'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.56%  '
